# "Removing "final" to VO getters and setters." -- database-comparison.xlsx
# (commit message refers to the source project; the actual spreadsheet edit
#  updates a couple of index-optimization feature cells and scrolls/selects
#  a different part of the sheet.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Content changes -------------------------------------------------
# Row 55 "Indexes on Expressions / Functional Indexes" -> MySQL (col G)
# used to read "Since 8.0", now gets the more precise "Since 8.0.15".
$ws.Range("A55").Value = "Indexes on Expressions / Functional Indexes"
$ws.Range("G55").Value = "Since 8.0.15"

# Row 58 "Indexes on Virtual Columns" -> PostgreSQL (col D) used to read
# "No", now reads "Since 12" (and is recoloured like the other "Yes"
# style cells, e.g. E58, instead of the "No" style).
$ws.Range("E58").Copy()
$ws.Range("D58").PasteSpecial(-4122)
$ws.Range("D58").Value = "Since 12"

# --- Column widths -----------------------------------------------------
# Column D (4) widened, and column G (7) is split off from the former
# G:H merged width band and widened on its own.
$ws.Columns.Item(4).ColumnWidth = 11.6
$ws.Columns.Item(7).ColumnWidth = 11.05

# --- View state ----------------------------------------------------
# Scroll position moved up and selection moved from A59 to A56.
$win = $excel.ActiveWindow
$win.ScrollRow = 46
$win.ScrollColumn = 1
[void]$ws.Range("A56").Select()
